$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 9.385801
$ws.Range("H2").Value = 18.771602
$ws.Range("I2").Value = 0.06848010139180623
$ws.Range("J2").Value = 0.04853166157087635
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 14.989415
$ws.Range("N2").Value = 29.97883
$ws.Range("O2").Value = 0.05547446260572933
$ws.Range("P2").Value = 0.03893791130463959
$ws.Range("Q2").Value = 140.687666296415
$ws.Range("R2").Value = 562.7506651856601
$ws.Range("S2").Value = 0.003798896823896307
$ws.Range("T2").Value = 0.001889721533713569

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 9.385801
$ws.Range("H3").Value = 18.771602
$ws.Range("I3").Value = 0.06848010139180623
$ws.Range("J3").Value = 0.04853166157087635
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 58.255493
$ws.Range("N3").Value = 174.766479
$ws.Range("O3").Value = 0.2155982850569436
$ws.Range("P3").Value = 0.2269949046819425
$ws.Range("Q3").Value = 546.7744644548931
$ws.Range("R3").Value = 3280.646786729358
$ws.Range("S3").Value = 0.01476419242059904
$ws.Range("T3").Value = 0.01101643989233737

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 9.385801
$ws.Range("H4").Value = 18.771602
$ws.Range("I4").Value = 0.06848010139180623
$ws.Range("J4").Value = 0.04853166157087635
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 80.178917
$ws.Range("N4").Value = 240.536751
$ws.Range("O4").Value = 0.2967348847759819
$ws.Range("P4").Value = 0.3124204205415681
$ws.Range("Q4").Value = 752.543359357517
$ws.Range("R4").Value = 4515.260156145103
$ws.Range("S4").Value = 0.02032043499594518
$ws.Range("T4").Value = 0.01516228211755425

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 9.385801
$ws.Range("H5").Value = 18.771602
$ws.Range("I5").Value = 0.06848010139180623
$ws.Range("J5").Value = 0.04853166157087635
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 79.34548433333335
$ws.Range("N5").Value = 238.036453
$ws.Range("O5").Value = 0.2936504262229702
$ws.Range("P5").Value = 0.3091729161606711
$ws.Range("Q5").Value = 744.7209262012846
$ws.Range("R5").Value = 4468.325557207707
$ws.Range("S5").Value = 0.02010921096149611
$ws.Range("T5").Value = 0.01500467533399062

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 9.385801
$ws.Range("H6").Value = 18.771602
$ws.Range("I6").Value = 0.06848010139180623
$ws.Range("J6").Value = 0.04853166157087635
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 11.726012
$ws.Range("N6").Value = 35.17803600000001
$ws.Range("O6").Value = 0.04339690469630293
$ws.Range("P6").Value = 0.04569088405516222
$ws.Range("Q6").Value = 110.058015155612
$ws.Range("R6").Value = 660.3480909336722
$ws.Range("S6").Value = 0.002971824433693376
$ws.Range("T6").Value = 0.002217454521839283

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 9.385801
$ws.Range("H7").Value = 18.771602
$ws.Range("I7").Value = 0.06848010139180623
$ws.Range("J7").Value = 0.04853166157087635
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 25.708558
$ws.Range("N7").Value = 51.417116
$ws.Range("O7").Value = 0.09514503664207198
$ws.Range("P7").Value = 0.0667829632560165
$ws.Range("Q7").Value = 241.295409384958
$ws.Range("R7").Value = 965.1816375398321
$ws.Range("S7").Value = 0.006515541756176208
$ws.Range("T7").Value = 0.003241088171441263

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 36.44655933333333
$ws.Range("H8").Value = 109.339678
$ws.Range("I8").Value = 0.2659191345021217
$ws.Range("J8").Value = 0.2826842508681249
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 14.989415
$ws.Range("N8").Value = 29.97883
$ws.Range("O8").Value = 0.05547446260572933
$ws.Range("P8").Value = 0.03893791130463959
$ws.Range("Q8").Value = 546.3126031694567
$ws.Range("R8").Value = 3277.875619016741
$ws.Range("S8").Value = 0.01475172108308586
$ws.Range("T8").Value = 0.01100713428752153

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 36.44655933333333
$ws.Range("H9").Value = 109.339678
$ws.Range("I9").Value = 0.2659191345021217
$ws.Range("J9").Value = 0.2826842508681249
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 58.255493
$ws.Range("N9").Value = 174.766479
$ws.Range("O9").Value = 0.2155982850569436
$ws.Range("P9").Value = 0.2269949046819425
$ws.Range("Q9").Value = 2123.212282117085
$ws.Range("R9").Value = 19108.91053905376
$ws.Range("S9").Value = 0.05733170936248416
$ws.Range("T9").Value = 0.06416788458089633

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 36.44655933333333
$ws.Range("H10").Value = 109.339678
$ws.Range("I10").Value = 0.2659191345021217
$ws.Range("J10").Value = 0.2826842508681249
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 80.178917
$ws.Range("N10").Value = 240.536751
$ws.Range("O10").Value = 0.2967348847759819
$ws.Range("P10").Value = 0.3124204205415681
$ws.Range("Q10").Value = 2922.245655722908
$ws.Range("R10").Value = 26300.21090150618
$ws.Range("S10").Value = 0.07890748373621592
$ws.Range("T10").Value = 0.08831633253669771

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 36.44655933333333
$ws.Range("H11").Value = 109.339678
$ws.Range("I11").Value = 0.2659191345021217
$ws.Range("J11").Value = 0.2826842508681249
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 79.34548433333335
$ws.Range("N11").Value = 238.036453
$ws.Range("O11").Value = 0.2936504262229702
$ws.Range("P11").Value = 0.3091729161606711
$ws.Range("Q11").Value = 2891.869902586904
$ws.Range("R11").Value = 26026.82912328214
$ws.Range("S11").Value = 0.07808726718739137
$ws.Range("T11").Value = 0.0873983141935929

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 36.44655933333333
$ws.Range("H12").Value = 109.339678
$ws.Range("I12").Value = 0.2659191345021217
$ws.Range("J12").Value = 0.2826842508681249
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 11.726012
$ws.Range("N12").Value = 35.17803600000001
$ws.Range("O12").Value = 0.04339690469630293
$ws.Range("P12").Value = 0.04569088405516222
$ws.Range("Q12").Value = 427.3727921013788
$ws.Range("R12").Value = 3846.355128912409
$ws.Range("S12").Value = 0.01154006733691193
$ws.Range("T12").Value = 0.01291609333063588

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 36.44655933333333
$ws.Range("H13").Value = 109.339678
$ws.Range("I13").Value = 0.2659191345021217
$ws.Range("J13").Value = 0.2826842508681249
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 25.708558
$ws.Range("N13").Value = 51.417116
$ws.Range("O13").Value = 0.09514503664207198
$ws.Range("P13").Value = 0.0667829632560165
$ws.Range("Q13").Value = 936.9884845214414
$ws.Range("R13").Value = 5621.930907128648
$ws.Range("S13").Value = 0.02530088579603243
$ws.Range("T13").Value = 0.01887849193878053

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 35.10635633333334
$ws.Range("H14").Value = 105.319069
$ws.Range("I14").Value = 0.2561408281726349
$ws.Range("J14").Value = 0.2722894622242564
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 14.989415
$ws.Range("N14").Value = 29.97883
$ws.Range("O14").Value = 0.05547446260572933
$ws.Range("P14").Value = 0.03893791130463959
$ws.Range("Q14").Value = 526.2237442182118
$ws.Range("R14").Value = 3157.342465309271
$ws.Range("S14").Value = 0.01420927479426338
$ws.Range("T14").Value = 0.01060238292927611

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 35.10635633333334
$ws.Range("H15").Value = 105.319069
$ws.Range("I15").Value = 0.2561408281726349
$ws.Range("J15").Value = 0.2722894622242564
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 58.255493
$ws.Range("N15").Value = 174.766479
$ws.Range("O15").Value = 0.2155982850569436
$ws.Range("P15").Value = 0.2269949046819425
$ws.Range("Q15").Value = 2045.138095632006
$ws.Range("R15").Value = 18406.24286068805
$ws.Range("S15").Value = 0.05522352328708537
$ws.Range("T15").Value = 0.06180832052349245

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 35.10635633333334
$ws.Range("H16").Value = 105.319069
$ws.Range("I16").Value = 0.2561408281726349
$ws.Range("J16").Value = 0.2722894622242564
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 80.178917
$ws.Range("N16").Value = 240.536751
$ws.Range("O16").Value = 0.2967348847759819
$ws.Range("P16").Value = 0.3124204205415681
$ws.Range("Q16").Value = 2814.789630622758
$ws.Range("R16").Value = 25333.10667560482
$ws.Range("S16").Value = 0.07600591913423141
$ws.Range("T16").Value = 0.0850687882971396

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 35.10635633333334
$ws.Range("H17").Value = 105.319069
$ws.Range("I17").Value = 0.2561408281726349
$ws.Range("J17").Value = 0.2722894622242564
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 79.34548433333335
$ws.Range("N17").Value = 238.036453
$ws.Range("O17").Value = 0.2936504262229702
$ws.Range("P17").Value = 0.3091729161606711
$ws.Range("Q17").Value = 2785.530846446918
$ws.Range("R17").Value = 25069.77761802226
$ws.Range("S17").Value = 0.07521586336599882
$ws.Range("T17").Value = 0.08418452707569424

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 35.10635633333334
$ws.Range("H18").Value = 105.319069
$ws.Range("I18").Value = 0.2561408281726349
$ws.Range("J18").Value = 0.2722894622242564
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 11.726012
$ws.Range("N18").Value = 35.17803600000001
$ws.Range("O18").Value = 0.04339690469630293
$ws.Range("P18").Value = 0.04569088405516222
$ws.Range("Q18").Value = 411.6575556409428
$ws.Range("R18").Value = 3704.918000768485
$ws.Range("S18").Value = 0.01111571910903994
$ws.Range("T18").Value = 0.01244114624793097

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 35.10635633333334
$ws.Range("H19").Value = 105.319069
$ws.Range("I19").Value = 0.2561408281726349
$ws.Range("J19").Value = 0.2722894622242564
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 25.708558
$ws.Range("N19").Value = 51.417116
$ws.Range("O19").Value = 0.09514503664207198
$ws.Range("P19").Value = 0.0667829632560165
$ws.Range("Q19").Value = 902.5337979641674
$ws.Range("R19").Value = 5415.202787785005
$ws.Range("S19").Value = 0.02437052848201602
$ws.Range("T19").Value = 0.01818429715072301

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 31.58644433333333
$ws.Range("H20").Value = 94.759333
$ws.Range("I20").Value = 0.2304590637020015
$ws.Range("J20").Value = 0.244988567296386
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 14.989415
$ws.Range("N20").Value = 29.97883
$ws.Range("O20").Value = 0.05547446260572933
$ws.Range("P20").Value = 0.03893791130463959
$ws.Range("Q20").Value = 473.4623224867317
$ws.Range("R20").Value = 2840.77393492039
$ws.Range("S20").Value = 0.01278459271148808
$ws.Range("T20").Value = 0.009539343104037408

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 31.58644433333333
$ws.Range("H21").Value = 94.759333
$ws.Range("I21").Value = 0.2304590637020015
$ws.Range("J21").Value = 0.244988567296386
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 58.255493
$ws.Range("N21").Value = 174.766479
$ws.Range("O21").Value = 0.2155982850569436
$ws.Range("P21").Value = 0.2269949046819425
$ws.Range("Q21").Value = 1840.08388675539
$ws.Range("R21").Value = 16560.75498079851
$ws.Range("S21").Value = 0.04968657890998045
$ws.Range("T21").Value = 0.0556111564816088

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 31.58644433333333
$ws.Range("H22").Value = 94.759333
$ws.Range("I22").Value = 0.2304590637020015
$ws.Range("J22").Value = 0.244988567296386
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 80.178917
$ws.Range("N22").Value = 240.536751
$ws.Range("O22").Value = 0.2967348847759819
$ws.Range("P22").Value = 0.3124204205415681
$ws.Range("Q22").Value = 2532.566898527454
$ws.Range("R22").Value = 22793.10208674708
$ws.Range("S22").Value = 0.0683852437131941
$ws.Range("T22").Value = 0.07653943122261318

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 31.58644433333333
$ws.Range("H23").Value = 94.759333
$ws.Range("I23").Value = 0.2304590637020015
$ws.Range("J23").Value = 0.244988567296386
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 79.34548433333335
$ws.Range("N23").Value = 238.036453
$ws.Range("O23").Value = 0.2936504262229702
$ws.Range("P23").Value = 0.3091729161606711
$ws.Range("Q23").Value = 2506.241723996206
$ws.Range("R23").Value = 22556.17551596585
$ws.Range("S23").Value = 0.06767440228303938
$ws.Range("T23").Value = 0.0757438297770485

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 31.58644433333333
$ws.Range("H24").Value = 94.759333
$ws.Range("I24").Value = 0.2304590637020015
$ws.Range("J24").Value = 0.244988567296386
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 11.726012
$ws.Range("N24").Value = 35.17803600000001
$ws.Range("O24").Value = 0.04339690469630293
$ws.Range("P24").Value = 0.04569088405516222
$ws.Range("Q24").Value = 370.3830252899987
$ws.Range("R24").Value = 3333.447227609989
$ws.Range("S24").Value = 0.01000121002387496
$ws.Range("T24").Value = 0.01119374422317948

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 31.58644433333333
$ws.Range("H25").Value = 94.759333
$ws.Range("I25").Value = 0.2304590637020015
$ws.Range("J25").Value = 0.244988567296386
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 25.708558
$ws.Range("N25").Value = 51.417116
$ws.Range("O25").Value = 0.09514503664207198
$ws.Range("P25").Value = 0.0667829632560165
$ws.Range("Q25").Value = 812.0419361572713
$ws.Range("R25").Value = 4872.251616943628
$ws.Range("S25").Value = 0.02192703606042454
$ws.Range("T25").Value = 0.01636106248789868

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 9.533863
$ws.Range("H26").Value = 28.601589
$ws.Range("I26").Value = 0.0695603822087843
$ws.Range("J26").Value = 0.07394588047079305
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 14.989415
$ws.Range("N26").Value = 29.97883
$ws.Range("O26").Value = 0.05547446260572933
$ws.Range("P26").Value = 0.03893791130463959
$ws.Range("Q26").Value = 142.907029060145
$ws.Range("R26").Value = 857.4421743608701
$ws.Range("S26").Value = 0.003858824821681445
$ws.Range("T26").Value = 0.00287929813511522

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 9.533863
$ws.Range("H27").Value = 28.601589
$ws.Range("I27").Value = 0.0695603822087843
$ws.Range("J27").Value = 0.07394588047079305
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 58.255493
$ws.Range("N27").Value = 174.766479
$ws.Range("O27").Value = 0.2155982850569436
$ws.Range("P27").Value = 0.2269949046819425
$ws.Range("Q27").Value = 555.399889259459
$ws.Range("R27").Value = 4998.599003335131
$ws.Range("S27").Value = 0.01499709911211943
$ws.Range("T27").Value = 0.01678533808908998

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 9.533863
$ws.Range("H28").Value = 28.601589
$ws.Range("I28").Value = 0.0695603822087843
$ws.Range("J28").Value = 0.07394588047079305
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 80.178917
$ws.Range("N28").Value = 240.536751
$ws.Range("O28").Value = 0.2967348847759819
$ws.Range("P28").Value = 0.3124204205415681
$ws.Range("Q28").Value = 764.414810166371
$ws.Range("R28").Value = 6879.733291497339
$ws.Range("S28").Value = 0.02064099199969687
$ws.Range("T28").Value = 0.02310220307400169

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 9.533863
$ws.Range("H29").Value = 28.601589
$ws.Range("I29").Value = 0.0695603822087843
$ws.Range("J29").Value = 0.07394588047079305
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 79.34548433333335
$ws.Range("N29").Value = 238.036453
$ws.Range("O29").Value = 0.2936504262229702
$ws.Range("P29").Value = 0.3091729161606711
$ws.Range("Q29").Value = 756.4689773026465
$ws.Range("R29").Value = 6808.220795723818
$ws.Range("S29").Value = 0.02042643588384222
$ws.Range("T29").Value = 0.02286206350322351

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 9.533863
$ws.Range("H30").Value = 28.601589
$ws.Range("I30").Value = 0.0695603822087843
$ws.Range("J30").Value = 0.07394588047079305
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 11.726012
$ws.Range("N30").Value = 35.17803600000001
$ws.Range("O30").Value = 0.04339690469630293
$ws.Range("P30").Value = 0.04569088405516222
$ws.Range("Q30").Value = 111.794191944356
$ws.Range("R30").Value = 1006.147727499204
$ws.Range("S30").Value = 0.003018705277353018
$ws.Range("T30").Value = 0.003378652650947889

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 9.533863
$ws.Range("H31").Value = 28.601589
$ws.Range("I31").Value = 0.0695603822087843
$ws.Range("J31").Value = 0.07394588047079305
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 25.708558
$ws.Range("N31").Value = 51.417116
$ws.Range("O31").Value = 0.09514503664207198
$ws.Range("P31").Value = 0.0667829632560165
$ws.Range("Q31").Value = 245.101869899554
$ws.Range("R31").Value = 1470.611219397324
$ws.Range("S31").Value = 0.006618325114091314
$ws.Range("T31").Value = 0.004938325018414761

$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 14.999783
$ws.Range("H32").Value = 29.999566
$ws.Range("I32").Value = 0.1094404900226514
$ws.Range("J32").Value = 0.07756017756956324
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 14.989415
$ws.Range("N32").Value = 29.97883
$ws.Range("O32").Value = 0.05547446260572933
$ws.Range("P32").Value = 0.03893791130463959
$ws.Range("Q32").Value = 224.837972296945
$ws.Range("R32").Value = 899.35188918778
$ws.Range("S32").Value = 0.006071152371314267
$ws.Range("T32").Value = 0.00302003131497575

$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 14.999783
$ws.Range("H33").Value = 29.999566
$ws.Range("I33").Value = 0.1094404900226514
$ws.Range("J33").Value = 0.07756017756956324
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 58.255493
$ws.Range("N33").Value = 174.766479
$ws.Range("O33").Value = 0.2155982850569436
$ws.Range("P33").Value = 0.2269949046819425
$ws.Range("Q33").Value = 873.819753558019
$ws.Range("R33").Value = 5242.918521348114
$ws.Range("S33").Value = 0.02359518196467518
$ws.Range("T33").Value = 0.01760576511451754

$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 14.999783
$ws.Range("H34").Value = 29.999566
$ws.Range("I34").Value = 0.1094404900226514
$ws.Range("J34").Value = 0.07756017756956324
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 80.178917
$ws.Range("N34").Value = 240.536751
$ws.Range("O34").Value = 0.2967348847759819
$ws.Range("P34").Value = 0.3124204205415681
$ws.Range("Q34").Value = 1202.666356175011
$ws.Range("R34").Value = 7215.998137050065
$ws.Range("S34").Value = 0.03247481119669846
$ws.Range("T34").Value = 0.02423138329356164

$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 14.999783
$ws.Range("H35").Value = 29.999566
$ws.Range("I35").Value = 0.1094404900226514
$ws.Range("J35").Value = 0.07756017756956324
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 79.34548433333335
$ws.Range("N35").Value = 238.036453
$ws.Range("O35").Value = 0.2936504262229702
$ws.Range("P35").Value = 0.3091729161606711
$ws.Range("Q35").Value = 1190.1650470299
$ws.Range("R35").Value = 7140.990282179398
$ws.Range("S35").Value = 0.03213724654120229
$ws.Range("T35").Value = 0.02397950627712134

$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 14.999783
$ws.Range("H36").Value = 29.999566
$ws.Range("I36").Value = 0.1094404900226514
$ws.Range("J36").Value = 0.07756017756956324
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 11.726012
$ws.Range("N36").Value = 35.17803600000001
$ws.Range("O36").Value = 0.04339690469630293
$ws.Range("P36").Value = 0.04569088405516222
$ws.Range("Q36").Value = 175.887635455396
$ws.Range("R36").Value = 1055.325812732376
$ws.Range("S36").Value = 0.004749378515429692
$ws.Range("T36").Value = 0.003543793080628707

$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 14.999783
$ws.Range("H37").Value = 29.999566
$ws.Range("I37").Value = 0.1094404900226514
$ws.Range("J37").Value = 0.07756017756956324
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 25.708558
$ws.Range("N37").Value = 51.417116
$ws.Range("O37").Value = 0.09514503664207198
$ws.Range("P37").Value = 0.0667829632560165
$ws.Range("Q37").Value = 385.622791242914
$ws.Range("R37").Value = 1542.491164971656
$ws.Range("S37").Value = 0.01041271943333148
$ws.Range("T37").Value = 0.005179698488758257
